# Section.xlsx: replace the "District Code / Section Name" sample data
# with an empty "Name / District ID" import template (see commit
# "Added Import For Church/District/DistrictState/Section/").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header row.
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "District ID"

# Drop the old sample rows (D001/Saltlake/Esplanade/Park Street) so the
# sheet is left with just the header - used range shrinks to A1:B1.
$ws.Range("A2:B4").ClearContents()

# Column A is no longer sized to fit the old "District Code" values;
# widen it to fit the new "Name" header instead.
$ws.Columns.Item(1).ColumnWidth = 14.7

# Active selection ends up on D8 after the edit.
[void]$ws.Range("D8").Select()
